$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used for PasteSpecial calls below.
$xlPasteFormats = -4122

# --- Copy the "Posts.js Model" / "Users.js Model" space-efficiency table from
# --- columns A:B into new columns D:E, to report the updated metric numbers.

# Row 1 - date header (copy date-header formatting from A1, set new date value)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D1").Value = 42704

# Row 2 - blank spacer row (copy formatting only)
$ws.Range("A2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial($xlPasteFormats) | Out-Null

# Row 3 - title (copy formatting + text)
$ws.Range("A3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D3").Value = "Space Effeciency"

# Row 5 - "Posts.js Model" label (no explicit style on source cell)
$ws.Range("D5").Value = "Posts.js Model"

# Row 6 - "each post takes" / "143 bytes of memory"
$ws.Range("D6").Value = "each post takes"
$ws.Range("B6").Copy() | Out-Null
$ws.Range("E6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E6").Value = "143 bytes of memory"

# Row 7 - blank formatted cell
$ws.Range("C7").Copy() | Out-Null
$ws.Range("E7").PasteSpecial($xlPasteFormats) | Out-Null

# Row 8 - blank formatted cells
$ws.Range("A8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B8").Copy() | Out-Null
$ws.Range("E8").PasteSpecial($xlPasteFormats) | Out-Null

# Row 9 - "Users.js Model" label / blank
$ws.Range("A9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D9").Value = "Users.js Model"
$ws.Range("B9").Copy() | Out-Null
$ws.Range("E9").PasteSpecial($xlPasteFormats) | Out-Null

# Row 10 - "each user takes" / updated memory metric
$ws.Range("A10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D10").Value = "each user takes"
$ws.Range("E10").Value = "117 bytes of memory"

# --- Widen / add columns to fit the new table ---
$ws.Columns.Item(2).ColumnWidth = 28
$ws.Columns.Item(4).ColumnWidth = 24.5
$ws.Columns.Item(5).ColumnWidth = 32.166666666666664

# --- Update the active selection to the new date cell ---
$ws.Range("D1").Select() | Out-Null
